$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 22, pushing existing rows 22-35 down to 23-36.
$ws.Rows.Item(22).Insert()

# Populate the newly inserted row 22 with the new data record.
$ws.Cells.Item(22, 1).Value = 4
$ws.Cells.Item(22, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(22, 3).Value = "Los Lagos"
$ws.Cells.Item(22, 4).Value = 44781
$ws.Cells.Item(22, 5).Value = 10
$ws.Cells.Item(22, 6).Value = 100112012
$ws.Cells.Item(22, 7).Value = "Espinaca"
$ws.Cells.Item(22, 8).Value = "Sin especificar"
$ws.Cells.Item(22, 9).Value = "Primera"
$ws.Cells.Item(22, 10).Value = 15
$ws.Cells.Item(22, 11).Value = 14000
$ws.Cells.Item(22, 12).Value = 14000
$ws.Cells.Item(22, 13).Value = 14000
$ws.Cells.Item(22, 14).Value = "$/cuna 10 kilos"
$ws.Cells.Item(22, 15).Value = "Región Metropolitana"
$ws.Cells.Item(22, 16).Value = 1400
$ws.Cells.Item(22, 17).Value = 10
$ws.Cells.Item(22, 18).Value = "Hortaliza"
